$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("MissingTableOfForAssertionTable")

# Replace the "Assert" label with "Then" to match BDD syntax
$ws.Range("A7").Value = "Then"

# Update the conditional formatting rule that matched on "Assert"
$fcs = $ws.Range("A1:XFD1048576").FormatConditions
for ($i = 1; $i -le $fcs.Count; $i++) {
    $fc = $fcs.Item($i)
    if ($fc.Formula1 -eq '="Assert"') {
        $fc.Formula1 = '="Then"'
    }
}
